$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 151; existing rows 151-166 shift down to 152-167.
$ws.Rows(151).Insert()

# Copy the date-cell number format from the row below (now row 152, formerly row 151)
# before writing the new date value so the new cell keeps the same "Fecha" style.
$ws.Range("D151").NumberFormat = $ws.Range("D152").NumberFormat

# Populate the new row 151 with the new weekly record.
$ws.Range("A151").Value2 = 7
$ws.Range("B151").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C151").Value2 = "Ñuble"
$ws.Range("D151").Value2 = 44449
$ws.Range("E151").Value2 = 16
$ws.Range("F151").Value2 = 100114013
$ws.Range("G151").Value2 = "Zanahoria"
$ws.Range("H151").Value2 = "Sin especificar"
$ws.Range("I151").Value2 = "Primera"
$ws.Range("J151").Value2 = 160
$ws.Range("K151").Value2 = 5000
$ws.Range("L151").Value2 = 5500
$ws.Range("M151").Value2 = 5250
$ws.Range("N151").Value2 = "$/saco 20 kilos"
$ws.Range("O151").Value2 = "Provincia de Diguillín"
$ws.Range("P151").Value2 = 262
$ws.Range("Q151").Value2 = 20
$ws.Range("R151").Value2 = "Hortaliza"
